# Insert a new data row at row 4 (pushing existing rows 4..102 down to 5..103)
# and populate the new row with the latest price entry, matching the
# "Fruta / hortaliza, semanal" update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 4 and below down by one row.
$ws.Rows("4:4").Insert()

# Populate the newly inserted row 4 with the new record.
$ws.Range("A4").Value = 10
$ws.Range("B4").Value = "Vega Modelo de Temuco"
$ws.Range("C4").Value = "La Araucanía"
$ws.Range("D4").Value = 44882
$ws.Range("E4").Value = 9
$ws.Range("F4").Value = 100112022
$ws.Range("G4").Value = "Arveja Verde"
$ws.Range("H4").Value = "Sin especificar"
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 80
$ws.Range("K4").Value = 16000
$ws.Range("L4").Value = 16000
$ws.Range("M4").Value = 16000
$ws.Range("N4").Value = "$/saco 25 kilos"
$ws.Range("O4").Value = "Región Metropolitana"
$ws.Range("P4").Value = 640
$ws.Range("Q4").Value = 25
$ws.Range("R4").Value = "Hortaliza"
